# Auto-generated edit script: updates market-price-derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, reflecting a refreshed
# data pull from the scheduled runner (chore: update Sheets via scheduled runner).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 280.81818
$ws.Range("J28").Value = 225
$ws.Range("L28").Value = 225
$ws.Range("N28").Value = -1195
$ws.Range("H32").Value = 1688.3334
$ws.Range("I32").Value = 2116.6667
$ws.Range("K32").Value = 2116.6667
$ws.Range("M32").Value = -1790.6667
$ws.Range("H40").Value = 2749.5
$ws.Range("J40").Value = 2666
$ws.Range("L40").Value = 2666
$ws.Range("N40").Value = -3016
$ws.Range("H106").Value = 3992.0476
$ws.Range("I106").Value = 3160.7334
$ws.Range("K106").Value = 3160.7334
$ws.Range("M106").Value = -2529.7334
$ws.Range("H118").Value = 591.25
$ws.Range("I118").Value = 591.25
$ws.Range("K118").Value = 1773.75
$ws.Range("M118").Value = -116.75
$ws.Range("H127").Value = 3347.1428
$ws.Range("I127").Value = 3738.3333
$ws.Range("J127").Value = 1000
$ws.Range("K127").Value = 11214.9999
$ws.Range("L127").Value = 3000
$ws.Range("M127").Value = -6254.999899999999
$ws.Range("N127").Value = -12920
$ws.Range("H129").Value = 1305.9286
$ws.Range("J129").Value = 1305.9286
$ws.Range("L129").Value = 3917.7858
$ws.Range("N129").Value = -13917.7858
$ws.Range("H137").Value = 2385.5715
$ws.Range("I137").Value = 1939.8
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 5819.4
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = -3269.4
$ws.Range("N137").Value = -15600
$ws.Range("H138").Value = 2908.6292
$ws.Range("I138").Value = 4007.5789
$ws.Range("K138").Value = 12022.7367
$ws.Range("M138").Value = -6882.736699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4125.72
$ws.Range("I32").Value = 4125.72
$ws.Range("K32").Value = 4125.72
$ws.Range("M32").Value = -3838.72
$ws.Range("H61").Value = 24120.756
$ws.Range("I61").Value = 36588.043
$ws.Range("J61").Value = 3638.7856
$ws.Range("K61").Value = 36588.043
$ws.Range("L61").Value = 3638.7856
$ws.Range("M61").Value = -36376.043
$ws.Range("N61").Value = -4062.7856
$ws.Range("H136").Value = 24120.756
$ws.Range("I136").Value = 36588.043
$ws.Range("J136").Value = 3638.7856
$ws.Range("K136").Value = 109764.129
$ws.Range("L136").Value = 10916.3568
$ws.Range("M136").Value = -107214.129
$ws.Range("N136").Value = -16016.3568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 3000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2887
$ws.Range("H76").Value = 22500
$ws.Range("J76").Value = 22500
$ws.Range("L76").Value = 22500
$ws.Range("N76").Value = -23130
$ws.Range("H79").Value = 22500
$ws.Range("J79").Value = 22500
$ws.Range("L79").Value = 22500
$ws.Range("N79").Value = -24684
$ws.Range("H82").Value = 14215.75
$ws.Range("I82").Value = 2954.3333
$ws.Range("K82").Value = 2954.3333
$ws.Range("M82").Value = -2571.3333
$ws.Range("H85").Value = 14215.75
$ws.Range("I85").Value = 2954.3333
$ws.Range("K85").Value = 2954.3333
$ws.Range("M85").Value = -1628.3333
$ws.Range("H86").Value = 667929
$ws.Range("I86").Value = 1858
$ws.Range("J86").Value = 1334000
$ws.Range("K86").Value = 1858
$ws.Range("L86").Value = 1334000
$ws.Range("M86").Value = -735
$ws.Range("N86").Value = -1336246
$ws.Range("H89").Value = 667929
$ws.Range("I89").Value = 1858
$ws.Range("J89").Value = 1334000
$ws.Range("K89").Value = 9290
$ws.Range("L89").Value = 6670000
$ws.Range("M89").Value = -3674
$ws.Range("N89").Value = -6681232
$ws.Range("H134").Value = 4665
$ws.Range("I134").Value = 3774.4285
$ws.Range("K134").Value = 11323.2855
$ws.Range("M134").Value = -8788.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 50.125
$ws.Range("I7").Value = 63.2
$ws.Range("J7").Value = 28.333334
$ws.Range("K7").Value = 63.2
$ws.Range("L7").Value = 28.333334
$ws.Range("M7").Value = 49.8
$ws.Range("N7").Value = -254.333334
$ws.Range("H12").Value = 450
$ws.Range("I12").Value = 450
$ws.Range("K12").Value = 450
$ws.Range("M12").Value = -280
$ws.Range("H22").Value = 1779.7
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1866.3334
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1866.3334
$ws.Range("M22").Value = -650
$ws.Range("N22").Value = -2566.3334
$ws.Range("H31").Value = 2451.3438
$ws.Range("I31").Value = 2225.6875
$ws.Range("J31").Value = 2677
$ws.Range("K31").Value = 2225.6875
$ws.Range("L31").Value = 2677
$ws.Range("M31").Value = -1930.6875
$ws.Range("N31").Value = -3267
$ws.Range("H34").Value = 2451.3438
$ws.Range("I34").Value = 2225.6875
$ws.Range("J34").Value = 2677
$ws.Range("K34").Value = 2225.6875
$ws.Range("L34").Value = 2677
$ws.Range("M34").Value = -2023.6875
$ws.Range("N34").Value = -3081
$ws.Range("H86").Value = 2354.3845
$ws.Range("I86").Value = 2171
$ws.Range("J86").Value = 2511.5715
$ws.Range("K86").Value = 2171
$ws.Range("L86").Value = 2511.5715
$ws.Range("M86").Value = -1048
$ws.Range("N86").Value = -4757.5715
$ws.Range("H89").Value = 2354.3845
$ws.Range("I89").Value = 2171
$ws.Range("J89").Value = 2511.5715
$ws.Range("K89").Value = 10855
$ws.Range("L89").Value = 12557.8575
$ws.Range("M89").Value = -5239
$ws.Range("N89").Value = -23789.8575
$ws.Range("H132").Value = 2218.5833
$ws.Range("I132").Value = 1701.4286
$ws.Range("K132").Value = 5104.2858
$ws.Range("M132").Value = -2574.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 838.3
$ws.Range("I68").Value = 697.5714
$ws.Range("J68").Value = 1166.6666
$ws.Range("K68").Value = 2092.7142
$ws.Range("L68").Value = 3499.9998
$ws.Range("M68").Value = -1281.7142
$ws.Range("N68").Value = -5121.9998
$ws.Range("H71").Value = 838.3
$ws.Range("I71").Value = 697.5714
$ws.Range("J71").Value = 1166.6666
$ws.Range("K71").Value = 6278.1426
$ws.Range("L71").Value = 10499.9994
$ws.Range("M71").Value = -2222.1426
$ws.Range("N71").Value = -18611.9994
$ws.Range("H131").Value = 7827655.5
$ws.Range("J131").Value = 16406.22
$ws.Range("L131").Value = 49218.66
$ws.Range("N131").Value = -59298.66
$ws.Range("H136").Value = 1521.6428
$ws.Range("I136").Value = 1521.6428
$ws.Range("K136").Value = 4564.928400000001
$ws.Range("M136").Value = 535.0715999999993
$ws.Range("H141").Value = 4086.5
$ws.Range("I141").Value = 4226.4443
$ws.Range("K141").Value = 12679.3329
$ws.Range("M141").Value = -7499.332900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5027.2
$ws.Range("I70").Value = 4574.75
$ws.Range("J70").Value = 5191.727
$ws.Range("K70").Value = 4574.75
$ws.Range("L70").Value = 5191.727
$ws.Range("M70").Value = -4304.75
$ws.Range("N70").Value = -5731.727
$ws.Range("H73").Value = 5027.2
$ws.Range("I73").Value = 4574.75
$ws.Range("J73").Value = 5191.727
$ws.Range("K73").Value = 4574.75
$ws.Range("L73").Value = 5191.727
$ws.Range("M73").Value = -3638.75
$ws.Range("N73").Value = -7063.727
$ws.Range("H132").Value = 1838729.9
$ws.Range("I132").Value = 2572541.8
$ws.Range("J132").Value = 4200
$ws.Range("K132").Value = 7717625.399999999
$ws.Range("L132").Value = 12600
$ws.Range("M132").Value = -7715095.399999999
$ws.Range("N132").Value = -17660

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 614.8333
$ws.Range("I9").Value = 172.5
$ws.Range("J9").Value = 1499.5
$ws.Range("K9").Value = 172.5
$ws.Range("L9").Value = 1499.5
$ws.Range("M9").Value = 51.5
$ws.Range("N9").Value = -1947.5
$ws.Range("H94").Value = 59999.5
$ws.Range("J94").Value = 59999.5
$ws.Range("L94").Value = 59999.5
$ws.Range("N94").Value = -61351.5
$ws.Range("H122").Value = 1762.2
$ws.Range("I122").Value = 1484.9412
$ws.Range("K122").Value = 4454.8236
$ws.Range("M122").Value = -2004.8236
$ws.Range("H136").Value = 1543
$ws.Range("I136").Value = 1056.5
$ws.Range("K136").Value = 3169.5
$ws.Range("M136").Value = -619.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 47097
$ws.Range("J70").Value = 47097
$ws.Range("L70").Value = 47097
$ws.Range("N70").Value = -47727
$ws.Range("H73").Value = 47097
$ws.Range("J73").Value = 47097
$ws.Range("L73").Value = 47097
$ws.Range("N73").Value = -49281
$ws.Range("H113").Value = 537.8333
$ws.Range("I113").Value = 478.5
$ws.Range("J113").Value = 620.9
$ws.Range("K113").Value = 1435.5
$ws.Range("L113").Value = 1862.7
$ws.Range("M113").Value = 734.5
$ws.Range("N113").Value = -6202.7
$ws.Range("H126").Value = 8456.9375
$ws.Range("I126").Value = 9371.615
$ws.Range("K126").Value = 28114.845
$ws.Range("M126").Value = -25644.845
$ws.Range("H136").Value = 17362678
$ws.Range("I136").Value = 27778846
$ws.Range("J136").Value = 2401.5833
$ws.Range("K136").Value = 83336538
$ws.Range("L136").Value = 7204.749899999999
$ws.Range("M136").Value = -83333988
$ws.Range("N136").Value = -12304.7499
$ws.Range("H141").Value = 65450.934
$ws.Range("J141").Value = 65983.14
$ws.Range("L141").Value = 65983.14
$ws.Range("N141").Value = -76343.14

